$d = $word.ActiveDocument

# 1. Merge the split runs of the "Création de semaine type..." bullet into a single run
#    by replacing the whole sentence with itself (Find/Replace collapses formatting-identical
#    runs into one run).
$d.Content.Find.Execute(
    "Création de semaine type avec plages de disponibilités (disponibilité récurrentes)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Création de semaine type avec plages de disponibilités (disponibilité récurrentes)", 2) | Out-Null

# 2. Remove the three sub-bullets that are no longer wanted: "isRecurant", the nested
#    "Si oui: fixe ou flottant..." bullet, and "isDone". Deleting each paragraph's range
#    (including its trailing paragraph mark) removes the bullet entirely and promotes the
#    remaining list items up, which also naturally restores the "Ajout de tâches dans une
#    semaine" bullet back to list level 0.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like "*isRecurant*" -or $text -like "*Si oui: fixe ou flottant*" -or $text -like "*isDone*") {
        $para.Range.Delete()
    }
}

# 3. Remove the leftover "_GoBack" bookmark at the end of the document.
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
